# Update the two worksheets ("NBR" and "BAR") with the fixed-workflow
# results: shift the Cutoff index (column B) by +4, replace the
# Reaction_number values (column C) with the re-run results, and drop
# the now-unused trailing rows (17-20), shrinking each sheet from
# A1:C20 down to A1:C16.

$wb = $excel.ActiveWorkbook

# New column B (Cutoff) and C (Reaction_number) values for rows 2..16
# of each sheet. Column A (index) is left untouched.
$nbrB = @(5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19)
$nbrC = @(78, 77, 77, 77, 80, 79, 78, 78, 78, 77, 77, 77, 77, 77, 76)

$barB = @(5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19)
$barC = @(575, 571, 573, 572, 566, 566, 565, 565, 565, 564, 563, 563, 563, 564, 563)

function Update-Sheet {
    param(
        [string]$SheetName,
        [array]$BValues,
        [array]$CValues
    )

    $ws = $wb.Worksheets.Item($SheetName)

    for ($i = 0; $i -lt $BValues.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $BValues[$i]
        $ws.Cells.Item($row, 3).Value = $CValues[$i]
    }

    # Remove the rows that are no longer part of the result set
    # (previously rows 17-20, now beyond the new A1:C16 range).
    $ws.Range("A17:C20").EntireRow.Delete()
}

Update-Sheet "NBR" $nbrB $nbrC
Update-Sheet "BAR" $barB $barC
